$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("5002093-18.2019.8.21.0067", "9000698-20.2019.8.21.0067", "Migrado"),
    @("5002094-03.2019.8.21.0067", "9001132-09.2019.8.21.0067", "Migrado"),
    @("5000997-02.2018.8.21.0067", "9000692-47.2018.8.21.0067", "Migrado"),
    @("5002039-18.2020.8.21.0067", "9000407-83.2020.8.21.0067", "Migrado"),
    @("5000998-84.2018.8.21.0067", "9000870-93.2018.8.21.0067", "Migrado"),
    @("5002095-85.2019.8.21.0067", "9000724-18.2019.8.21.0067", "Migrado"),
    @("5002096-70.2019.8.21.0067", "9000532-85.2019.8.21.0067", "Migrado"),
    @("5002097-55.2019.8.21.0067", "9000464-38.2019.8.21.0067", "Migrado"),
    @("5053427-32.2018.8.21.0001", "9068801-20.2018.8.21.0001", "Migrado"),
    @("5035085-12.2014.8.21.0001", "0395366-43.2014.8.21.0001", "Digitalizado"),
    @("5013439-53.2008.8.21.0001", "1253081-20.2008.8.21.0001", "Digitalizado"),
    @("5013439-53.2008.8.21.0001", "1253081-20.2008.8.21.0001", "Digitalizado"),
    @("5013439-53.2008.8.21.0001", "1253081-20.2008.8.21.0001", "Digitalizado"),
    @("5013439-53.2008.8.21.0001", "1253081-20.2008.8.21.0001", "Digitalizado"),
    @("5013439-53.2008.8.21.0001", "1253081-20.2008.8.21.0001", "Digitalizado"),
    @("5014763-10.2010.8.21.0001", "3103141-73.2010.8.21.0001", "Digitalizado"),
    @("5014763-10.2010.8.21.0001", "3103141-73.2010.8.21.0001", "Digitalizado"),
    @("5014763-10.2010.8.21.0001", "3103141-73.2010.8.21.0001", "Digitalizado"),
    @("5029278-74.2015.8.21.0001", "0199864-35.2015.8.21.0001", "Digitalizado"),
    @("5002539-45.2007.8.21.0001", "0183341-26.2007.8.21.0001", "Digitalizado"),
    @("5020879-61.2012.8.21.0001", "0026384-42.2012.8.21.3001", "Digitalizado")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
